$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00158971523361902
$ws.Range("C2").Value = 0.0000691180536356096
$ws.Range("G2").Value = 0.999792645839093
$ws.Range("H2").Value = 0.44297760575062201
$ws.Range("I2").Value = 0.99930881946364403
$ws.Range("J2").Value = 0.000207354160906829
$ws.Range("L2").Value = 0.99647497926458395
$ws.Range("M2").Value = 0.0000691180536356096
$ws.Range("N2").Value = 0.000552944429084877
$ws.Range("O2").Value = 0.000276472214542438
$ws.Range("P2").Value = 0.000207354160906829
$ws.Range("Q2").Value = 0.000483826375449267
$ws.Range("R2").Value = 0.000207354160906829
$ws.Range("S2").Value = 0.000898534697262925
$ws.Range("T2").Value = 0.000345590268178048
$ws.Range("U2").Value = 0.000276472214542438
$ws.Range("V2").Value = 0.99640586121094799
$ws.Range("W2").Value = 0.91588332872546296
$ws.Range("X2").Value = 0.000967652750898535
$ws.Range("B3").Value = 0.997649986176389
$ws.Range("C3").Value = 0.99972352778545803
$ws.Range("D3").Value = 0.99923970141000795
$ws.Range("E3").Value = 0.000276472214542438
$ws.Range("F3").Value = 0.999792645839093
$ws.Range("H3").Value = 0.0000691180536356096
$ws.Range("I3").Value = 0.000276472214542438
$ws.Range("J3").Value = 0.000483826375449267
$ws.Range("K3").Value = 0.99149847940281999
$ws.Range("L3").Value = 0.00138236107271219
$ws.Range("M3").Value = 0.99875587503455898
$ws.Range("N3").Value = 0.000276472214542438
$ws.Range("O3").Value = 0.99896322919546598
$ws.Range("P3").Value = 0.000760298589991706
$ws.Range("Q3").Value = 0.000829416643627315
$ws.Range("R3").Value = 0.000552944429084877
$ws.Range("T3").Value = 0.000138236107271219
$ws.Range("U3").Value = 0.99903234724910195
$ws.Range("V3").Value = 0.000483826375449267
$ws.Range("W3").Value = 0.000138236107271219
$ws.Range("X3").Value = 0.99868675698092302
$ws.Range("B4").Value = 0.0000691180536356096
$ws.Range("F4").Value = 0.0000691180536356096
$ws.Range("G4").Value = 0.0000691180536356096
$ws.Range("H4").Value = 0.54810616533038403
$ws.Range("I4").Value = 0.000414708321813658
$ws.Range("J4").Value = 0.000345590268178048
$ws.Range("L4").Value = 0.000967652750898535
$ws.Range("M4").Value = 0.000552944429084877
$ws.Range("O4").Value = 0.0000691180536356096
$ws.Range("Q4").Value = 0.0000691180536356096
$ws.Range("R4").Value = 0.000483826375449267
$ws.Range("S4").Value = 0.99903234724910195
$ws.Range("T4").Value = 0.99944705557091496
$ws.Range("U4").Value = 0.000345590268178048
$ws.Range("V4").Value = 0.00304119435996682
$ws.Range("W4").Value = 0.080729886646392
$ws.Range("X4").Value = 0.000138236107271219
$ws.Range("B5").Value = 0.000414708321813658
$ws.Range("C5").Value = 0.000138236107271219
$ws.Range("D5").Value = 0.000552944429084877
$ws.Range("E5").Value = 0.99951617362455103
$ws.Range("F5").Value = 0.000138236107271219
$ws.Range("G5").Value = 0.0000691180536356096
$ws.Range("H5").Value = 0.000138236107271219
$ws.Range("J5").Value = 0.99875587503455898
$ws.Range("K5").Value = 0.00836328448990876
$ws.Range("L5").Value = 0.000829416643627315
$ws.Range("M5").Value = 0.000207354160906829
$ws.Range("N5").Value = 0.99910146530273702
$ws.Range("O5").Value = 0.000207354160906829
$ws.Range("P5").Value = 0.99896322919546598
$ws.Range("Q5").Value = 0.99854852087365198
$ws.Range("R5").Value = 0.99854852087365198
$ws.Range("S5").Value = 0.0000691180536356096
$ws.Range("T5").Value = 0.0000691180536356096
$ws.Range("U5").Value = 0.000207354160906829
$ws.Range("W5").Value = 0.0000691180536356096
